# tl_p028r.docx edit:
#  1. Add a comment (id 0) anchored on the "<" that opens "<caption><man>"
#     in the run just after the <link> hyperlink paragraph.
#  2. Fix up the sentence " a piece marked with that caliber weighs." ->
#     " the piece that carries the marked caliber weighs."
#
# (The hyperlink's relationship id shifting from rId6 to rId7 in the
# original commit is purely a side effect of Word inserting the new
# word/comments.xml relationship ahead of the hyperlink relationship when
# it renumbers the package on save -- it is not a content change we drive
# through the object model.)

$word.UserName = "Tianna Uchacz"
$word.UserInitials = "TU"

$d = $word.ActiveDocument

# --- 1. anchor + insert the comment ---------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("<caption><man>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not locate '<caption><man>' run"
}

# Comment range = just the leading "<" character of that run, matching
# the commentRangeStart/commentRangeEnd placement in the target edit.
$commentAnchor = $d.Range($rng.Start, $rng.Start + 1)

$commentText = "+catapanoth@gmail.com  We are translating here, but this caption is not rendering. We're not touching this mark-up (same across all versions), since it looks like you and Nick are working on figure/text/caption rendering." + [char]13 + "_Assigned to Terry Catapano_"

$d.Comments.Add($commentAnchor, $commentText) | Out-Null

# --- 2. caption sentence rewrite --------------------------------------
$d.Content.Find.Execute(" a piece marked with that caliber weighs.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " the piece that carries the marked caliber weighs.", 2) | Out-Null
